$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.027388102864774
$ws.Range("D2").Value = 1.03079118328872
$ws.Range("E2").Value = 1.027386650394401
$ws.Range("F2").Value = 1.025920811995852
$ws.Range("I2").Value = 1.029688338125137
$ws.Range("J2").Value = 1.032546254803232
$ws.Range("K2").Value = 1.033601172550376
$ws.Range("L2").Value = 1.030206529170885
$ws.Range("M2").Value = 1.028744970734982
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.029210797770532
$ws.Range("D3").Value = 1.032161511241162
$ws.Range("E3").Value = 1.028964918901997
$ws.Range("F3").Value = 1.028391611944887
$ws.Range("I3").Value = 1.030078471264797
$ws.Range("J3").Value = 1.034004631635953
$ws.Range("K3").Value = 1.034778421433916
$ws.Range("L3").Value = 1.031590430017403
$ws.Range("M3").Value = 1.031018671749184
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.030385608663026
$ws.Range("D4").Value = 1.033044086593452
$ws.Range("E4").Value = 1.029982246430839
$ws.Range("F4").Value = 1.029984943621277
$ws.Range("I4").Value = 1.030327797982162
$ws.Range("J4").Value = 1.034943526352196
$ws.Range("K4").Value = 1.035535549548722
$ws.Range("L4").Value = 1.032481525277856
$ws.Range("M4").Value = 1.032484215561182
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.030878422396736
$ws.Range("D5").Value = 1.033414151611036
$ws.Range("E5").Value = 1.030409011911067
$ws.Range("F5").Value = 1.030653514811744
$ws.Range("I5").Value = 1.030431875041359
$ws.Range("J5").Value = 1.035337114287145
$ws.Range("K5").Value = 1.035852753030572
$ws.Range("L5").Value = 1.032855110586847
$ws.Range("M5").Value = 1.033099001739449
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.030961105517148
$ws.Range("D6").Value = 1.033476230764102
$ws.Range("E6").Value = 1.030480614366127
$ws.Range("F6").Value = 1.030765697720659
$ws.Range("I6").Value = 1.030449306821587
$ws.Range("J6").Value = 1.03540313410845
$ws.Range("K6").Value = 1.035905949230033
$ws.Range("L6").Value = 1.032917777247752
$ws.Range("M6").Value = 1.033202150135097
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.030392197863891
$ws.Range("D7").Value = 1.033049035208335
$ws.Range("E7").Value = 1.029987952477596
$ws.Range("F7").Value = 1.029993882024683
$ws.Range("I7").Value = 1.030329191564081
$ws.Range("J7").Value = 1.034948789881744
$ws.Range("K7").Value = 1.035539792311753
$ws.Range("L7").Value = 1.032486521170753
$ws.Range("M7").Value = 1.032492435537804
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.028005056527365
$ws.Range("D8").Value = 1.031255154672267
$ws.Range("E8").Value = 1.027920856266567
$ws.Range("F8").Value = 1.02675697953702
$ws.Range("I8").Value = 1.029820834150756
$ws.Range("J8").Value = 1.033040119896044
$ws.Range("K8").Value = 1.033999997460164
$ws.Range("L8").Value = 1.030675142812246
$ws.Range("M8").Value = 1.029514578218143
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.023762369543413
$ws.Range("D9").Value = 1.028061826352144
$ws.Range("E9").Value = 1.024247507504275
$ws.Range("F9").Value = 1.021009788581188
$ws.Range("I9").Value = 1.028900911511854
$ws.Range("J9").Value = 1.029639378597748
$ws.Range("K9").Value = 1.03125051879773
$ws.Range("L9").Value = 1.027448906160816
$ws.Range("M9").Value = 1.024222054496853
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.020908066339006
$ws.Range("D10").Value = 1.02591019285647
$ws.Range("E10").Value = 1.021776651489312
$ws.Range("F10").Value = 1.017146806851837
$ws.Range("I10").Value = 1.028271021615735
$ws.Range("J10").Value = 1.02734583275971
$ws.Range("K10").Value = 1.029392232751239
$ws.Range("L10").Value = 1.025273850195779
$ws.Range("M10").Value = 1.020661146287537
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.019665658498755
$ws.Range("D11").Value = 1.024972881107236
$ws.Range("E11").Value = 1.020701263889504
$ws.Range("F11").Value = 1.015466083675844
$ws.Range("I11").Value = 1.027994246508848
$ws.Range("J11").Value = 1.026346168085579
$ws.Range("K11").Value = 1.028581351257984
$ws.Range("L11").Value = 1.024326024664469
$ws.Range("M11").Value = 1.019111017244328
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.019203171090338
$ws.Range("D12").Value = 1.024623853929373
$ws.Range("E12").Value = 1.020300968517793
$ws.Range("F12").Value = 1.014840536117971
$ws.Range("I12").Value = 1.027890827077006
$ws.Range("J12").Value = 1.025973840933575
$ws.Range("K12").Value = 1.028279198215507
$ws.Range("L12").Value = 1.023973034662332
$ws.Range("M12").Value = 1.01853395063955
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.019302422012985
$ws.Range("D13").Value = 1.024698761048725
$ws.Range("E13").Value = 1.020386872015344
$ws.Range("F13").Value = 1.014974775604908
$ws.Range("I13").Value = 1.02791303875128
$ws.Range("J13").Value = 1.026053752338059
$ws.Range("K13").Value = 1.028344054623748
$ws.Range("L13").Value = 1.024048794452074
$ws.Range("M13").Value = 1.018657791988022
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.019627449751082
$ws.Range("D14").Value = 1.024944048235637
$ws.Range("E14").Value = 1.020668192797182
$ws.Range("F14").Value = 1.015414401533084
$ws.Range("I14").Value = 1.027985710372296
$ws.Range("J14").Value = 1.026315412093499
$ws.Range("K14").Value = 1.028556394799374
$ws.Range("L14").Value = 1.024296865408213
$ws.Range("M14").Value = 1.019063343093954
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.019827576443929
$ws.Range("D15").Value = 1.025095062148494
$ws.Range("E15").Value = 1.0208414106628
$ws.Range("F15").Value = 1.015685102260009
$ws.Range("I15").Value = 1.02803040432104
$ws.Range("J15").Value = 1.026476495118641
$ws.Range("K15").Value = 1.028687097385053
$ws.Range("L15").Value = 1.024449586817804
$ws.Range("M15").Value = 1.019313045653444
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.020990381962374
$ws.Range("D16").Value = 1.025972278576817
$ws.Range("E16").Value = 1.021847903733155
$ws.Range("F16").Value = 1.017258177874437
$ws.Range("I16").Value = 1.028289304745706
$ws.Range("J16").Value = 1.027412037294556
$ws.Range("K16").Value = 1.029445915352658
$ws.Range("L16").Value = 1.025336625705897
$ws.Range("M16").Value = 1.020763846161115
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.021718025145026
$ws.Range("D17").Value = 1.026521008864444
$ws.Range("E17").Value = 1.022477763685838
$ws.Range("F17").Value = 1.018242745309833
$ws.Range("I17").Value = 1.028450622256203
$ws.Range("J17").Value = 1.027997109307553
$ws.Range("K17").Value = 1.029920219937332
$ws.Range("L17").Value = 1.025891416660287
$ws.Range("M17").Value = 1.021671660950311
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.022141824266523
$ws.Range("D18").Value = 1.026840531156086
$ws.Range("E18").Value = 1.022844622188338
$ws.Range("F18").Value = 1.018816255548458
$ws.Range("I18").Value = 1.028544327845645
$ws.Range("J18").Value = 1.028337742135138
$ws.Range("K18").Value = 1.0301962738641
$ws.Range("L18").Value = 1.026214438288225
$ws.Range("M18").Value = 1.02220038259124
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.022286224039145
$ws.Range("D19").Value = 1.026949388546213
$ws.Range("E19").Value = 1.022969622582725
$ws.Range("F19").Value = 1.019011678282565
$ws.Range("I19").Value = 1.028576213414405
$ws.Range("J19").Value = 1.028453783025156
$ws.Range("K19").Value = 1.030290299944441
$ws.Range("L19").Value = 1.026324482849074
$ws.Range("M19").Value = 1.02238053006396
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.021640020548739
$ws.Range("D20").Value = 1.02646219158951
$ws.Range("E20").Value = 1.022410240473024
$ws.Range("F20").Value = 1.018137190652185
$ws.Range("I20").Value = 1.028433354620343
$ws.Range("J20").Value = 1.027934401931179
$ws.Range("K20").Value = 1.029869393735215
$ws.Range("L20").Value = 1.025831952842935
$ws.Range("M20").Value = 1.021574343177671
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.019531765072958
$ws.Range("D21").Value = 1.024871841396066
$ws.Range("E21").Value = 1.020585374454092
$ws.Range("F21").Value = 1.015284977594884
$ws.Range("I21").Value = 1.027964327369926
$ws.Range("J21").Value = 1.026238387815021
$ws.Range("K21").Value = 1.02849389244507
$ws.Range("L21").Value = 1.024223840420023
$ws.Range("M21").Value = 1.018943954073924
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.01820041029841
$ws.Range("D22").Value = 1.023866893540325
$ws.Range("E22").Value = 1.019433087404771
$ws.Range("F22").Value = 1.013484410090491
$ws.Range("I22").Value = 1.027665881804098
$ws.Range("J22").Value = 1.02516619700473
$ws.Range("K22").Value = 1.02762352183618
$ws.Range("L22").Value = 1.023207390995284
$ws.Range("M22").Value = 1.017282698574559
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.018906747236364
$ws.Range("D23").Value = 1.024400119274367
$ws.Range("E23").Value = 1.020044411050158
$ws.Range("F23").Value = 1.014439629552488
$ws.Range("I23").Value = 1.027824432463408
$ws.Range("J23").Value = 1.025735147264079
$ws.Range("K23").Value = 1.02808545328471
$ws.Range("L23").Value = 1.023746746136102
$ws.Range("M23").Value = 1.018164079893409
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.021675269391704
$ws.Range("D24").Value = 1.026488770260924
$ws.Range("E24").Value = 1.022440752937626
$ws.Range("F24").Value = 1.018184888639803
$ws.Range("I24").Value = 1.028441158320963
$ws.Range("J24").Value = 1.027962738640553
$ws.Range("K24").Value = 1.029892361761262
$ws.Range("L24").Value = 1.025858823769387
$ws.Range("M24").Value = 1.021618319333701
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.024863651307452
$ws.Range("D25").Value = 1.028891310207258
$ws.Range("E25").Value = 1.025200936565805
$ws.Range("F25").Value = 1.022500954845393
$ws.Range("I25").Value = 1.029141632367965
$ws.Range("J25").Value = 1.030523111253799
$ws.Range("K25").Value = 1.031965711012572
$ws.Range("L25").Value = 1.028287152545051
$ws.Range("M25").Value = 1.025595874297711
